$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'284.15"
$ws.Range("E2").Value = "'2.54%"
$ws.Range("D3").Value = "'28.29"
$ws.Range("E3").Value = "'3.62%"
$ws.Range("D4").Value = "'5.056"
$ws.Range("E4").Value = "'4.45%"
$ws.Range("D5").Value = "'0.06645"
$ws.Range("E5").Value = "'5.07%"
$ws.Range("D6").Value = "'7.315"
$ws.Range("E6").Value = "'4.15%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.387"
$ws.Range("E7").Value = "'4.84%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9367"
$ws.Range("E8").Value = "'5.45%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1560"
$ws.Range("E9").Value = "'3.31%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.05969"
$ws.Range("E10").Value = "'10.67%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07621"
$ws.Range("E11").Value = "'2.42%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02902"
$ws.Range("E12").Value = "'0.01%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.08971"
$ws.Range("E13").Value = "'0.20%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001602"
$ws.Range("E14").Value = "'1.97%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04448"
$ws.Range("E15").Value = "'1.82%"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D16").Value = "'0.0006791"
$ws.Range("E16").Value = "'7.03%"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006122"
$ws.Range("E17").Value = "'1.48%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.441"
$ws.Range("E18").Value = "'-0.92%"
$ws.Range("B19").Value = "GateToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D19").Value = "'3.371"
$ws.Range("E19").Value = "'2.22%"
$ws.Range("D20").Value = "'2.237"
$ws.Range("E20").Value = "'0.16%"
$ws.Range("D21").Value = "'0.3215"
$ws.Range("E21").Value = "'1.40%"
$ws.Range("D22").Value = "'0.1274"
$ws.Range("E22").Value = "'-4.92%"
$ws.Range("D23").Value = "'4.029"
$ws.Range("E23").Value = "'3.09%"
$ws.Range("D24").Value = "'0.1520"
$ws.Range("E24").Value = "'0.87%"
$ws.Range("D25").Value = "'0.001170"
$ws.Range("E25").Value = "'-0.43%"
$ws.Range("D26").Value = "'0.004432"
$ws.Range("E26").Value = "'4.33%"
$ws.Range("D27").Value = "'0.0001231"
$ws.Range("E27").Value = "'4.34%"
$ws.Range("D28").Value = "'0.0001609"
$ws.Range("E28").Value = "'-2.41%"
$ws.Range("D40").Value = "'0.04171"
$ws.Range("E40").Value = "'4.40%"
$ws.Range("D41").Value = "'0.006220"
$ws.Range("E41").Value = "'-6.48%"
$ws.Range("D42").Value = "'0.1228"
$ws.Range("E42").Value = "'-12.00%"
$ws.Range("D43").Value = "'0.001988"
$ws.Range("E43").Value = "'-7.47%"
$ws.Range("D44").Value = "'0.01203"
$ws.Range("E44").Value = "'1.97%"
$ws.Range("D45").Value = "'0.00005492"
$ws.Range("E45").Value = "'-1.34%"
$ws.Range("D47").Value = "'0.01298"
$ws.Range("E47").Value = "'-29.84%"
